$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.890.28'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = '3.500.26'
$ws.Range('E3').Value = '  -1.93%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '603.16'
$ws.Range('E5').Value = '  -1.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '197.94'
$ws.Range('E6').Value = '  +6.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.625'
$ws.Range('E7').Value = '  +1.11%  '
$ws.Range('E9').Value = '  -2.71%  '
$ws.Range('E10').Value = '  +1.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '54.29'
$ws.Range('E11').Value = '  +0.56%  '
$ws.Range('E12').Value = '  -2.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.59'
$ws.Range('E13').Value = '  +0.91%  '
$ws.Range('D14').Value = '4.055.92'
$ws.Range('E14').Value = '  -1.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '595.62'
$ws.Range('E15').Value = '  +2.69%  '
$ws.Range('D16').Value = '69.967.05'
$ws.Range('E16').Value = '  -0.58%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.06'
$ws.Range('E17').Value = '  +0.34%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.62'
$ws.Range('E18').Value = '  -0.83%  '
$ws.Range('D19').Value = '3.501.10'
$ws.Range('E19').Value = '  -2.83%  '
$ws.Range('E20').Value = '  +0.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.990'
$ws.Range('E21').Value = '  -0.61%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '18.34'
$ws.Range('E22').Value = '  +5.66%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '104.61'
$ws.Range('E23').Value = '  +10.37%  '
$ws.Range('E24').Value = '  -3.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.03'
$ws.Range('E25').Value = '  +3.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.11'
$ws.Range('E26').Value = '  +4.88%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.99'
$ws.Range('E27').Value = '  +0.21%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.75'
$ws.Range('E28').Value = '  +3.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.63'
$ws.Range('E29').Value = '  +4.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.53'
$ws.Range('E30').Value = '  +22.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.25'
$ws.Range('E31').Value = '  +2.75%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.73'
$ws.Range('E32').Value = '  +3.65%  '
$ws.Range('E33').Value = '  +1.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.62'
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('D35').Value = '3.732.99'
$ws.Range('E35').Value = '  +5.85%  '
$ws.Range('D36').Value = '0.0₃0809'
$ws.Range('E36').Value = '  +2.93%  '
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '510.59'
$ws.Range('E38').Value = '  -4.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.392'
$ws.Range('E39').Value = '  -2.98%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.99'
$ws.Range('E40').Value = '  -7.99%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.72'
$ws.Range('E41').Value = '  -1.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.52'
$ws.Range('E42').Value = '  +0.06%  '
$ws.Range('E43').Value = '  +0.90%  '
$ws.Range('E44').Value = '  -0.48%  '
$ws.Range('E45').Value = '  -3.38%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.32'
$ws.Range('E46').Value = '  -3.49%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.140'
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('B48').Value = 'FirstDigitalUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.00'
$ws.Range('E48').Value = '  +0.24%  '
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.76'
$ws.Range('E49').Value = '  -5.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.02'
$ws.Range('E50').Value = '  -3.24%  '
$ws.Range('E51').Value = '  -1.94%  '
